# Add semantic_label metadata for the "tooth" rows (5-8) and move the
# active selection, per commit "add metadata to spreadsheet".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = "tooth 1"
$ws.Range("E6").Value = "tooth 2"
$ws.Range("E7").Value = "tooth 3"
$ws.Range("E8").Value = "tooth 4"

$ws.Range("A9").Select()
